$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (Price column D, Volume(1h) column E)
# Values that look numeric are written as text (NumberFormat "@") so Excel
# does not silently convert them into actual numbers, matching the source
# data which stores these as plain text strings.

# Row 2
$ws.Range("D2").Value = "26.896.92"
$ws.Range("E2").Value = "  -2.02%  "

# Row 3
$ws.Range("D3").Value = "1.833.19"
$ws.Range("E3").Value = "  -1.70%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.53"
$ws.Range("E5").Value = "  -1.48%  "

# Row 6
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4629"
$ws.Range("E7").Value = "  -0.64%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3672"
$ws.Range("E8").Value = "  -1.45%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07169"
$ws.Range("E9").Value = "  -2.74%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8770"
$ws.Range("E10").Value = "  -1.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07911"
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.61"
$ws.Range("E12").Value = "  -1.70%  "

# Row 13
$ws.Range("D13").Value = "1.860.30"
$ws.Range("E13").Value = "  -0.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("E14").Value = "  -1.31%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.397"
$ws.Range("E15").Value = "  -3.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.94"
$ws.Range("E16").Value = "  -5.18%  "

# Row 17
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008739"
$ws.Range("E18").Value = "  -1.93%  "

# Row 19
$ws.Range("E19").Value = "  +0.20%  "

# Row 20
$ws.Range("D20").Value = "26.917.43"
$ws.Range("E20").Value = "  -2.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").Value = "  -2.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.006"
$ws.Range("E22").Value = "  -2.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.44"
$ws.Range("E23").Value = "  -1.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.991"
$ws.Range("E24").Value = "  +4.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.90"
$ws.Range("E25").Value = "  -1.65%  "

# Row 26
$ws.Range("E26").Value = "  -1.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.970"
$ws.Range("E27").Value = "  -5.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.71"
$ws.Range("E28").Value = "  -2.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.953"
$ws.Range("E29").Value = "  -3.92%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08851"
$ws.Range("E30").Value = "  -0.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.127"
$ws.Range("E31").Value = "  +3.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7574"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.458"
$ws.Range("E33").Value = "  -0.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.128"
$ws.Range("E34").Value = "  -3.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.581"
$ws.Range("E35").Value = "  -2.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.086"
$ws.Range("E36").Value = "  +0.31%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01936"
$ws.Range("E37").Value = "  -1.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.925"
$ws.Range("E38").Value = "  -2.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05135"
$ws.Range("E39").Value = "  -2.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.908"
$ws.Range("E40").Value = "  -3.69%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4969"
$ws.Range("E41").Value = "  -4.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1599"
$ws.Range("E42").Value = "  -2.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.322"
$ws.Range("E43").Value = "  -0.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4684"
$ws.Range("E44").Value = "  -3.73%  "

# Row 45
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.07"
$ws.Range("E46").Value = "  -2.70%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.46"
$ws.Range("E47").Value = "  -1.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.614"
$ws.Range("E48").Value = "  -2.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06097"
$ws.Range("E49").Value = "  -2.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.69"
$ws.Range("E50").Value = "  -1.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.37"
$ws.Range("E51").Value = "  -2.37%  "
